$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 656, pushing existing rows 656:696 down to 657:697
$ws.Rows(656).Insert()

# Populate the newly inserted row 656 with the new price record
$ws.Cells.Item(656, 1).Value = 8
$ws.Cells.Item(656, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(656, 3).Value = "Coquimbo"
$ws.Cells.Item(656, 4).Value = 44585
$ws.Cells.Item(656, 5).Value = 4
$ws.Cells.Item(656, 6).Value = "Fruta"
$ws.Cells.Item(656, 7).Value = 100102
$ws.Cells.Item(656, 8).Value = "Cítricos"
$ws.Cells.Item(656, 9).Value = 100102005
$ws.Cells.Item(656, 10).Value = "Naranja"
$ws.Cells.Item(656, 11).Value = "Navel Late"
$ws.Cells.Item(656, 12).Value = "Primera"
$ws.Cells.Item(656, 13).Value = 16
$ws.Cells.Item(656, 14).Value = 235000
$ws.Cells.Item(656, 15).Value = 240000
$ws.Cells.Item(656, 16).Value = 237500
$ws.Cells.Item(656, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(656, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(656, 19).Value = 594
$ws.Cells.Item(656, 20).Value = 400
